# Datorama Creative Hierarchical workbook update
# "Added the display level changes to the framework"
#
# Renames the destination-table display labels for the raw delivery/
# conversion metrics so they're suffixed with "(IOne)" to disambiguate
# the source system, across every mapping sheet in the workbook.

$wb = $excel.ActiveWorkbook

# sheet index (1-based) -> list of (cell, newValue) edits to apply
$deliveryEdits = @(
    @{ Cell = "B3"; Value = "Impressions (IOne)" },
    @{ Cell = "B4"; Value = "Media Cost (IOne)" },
    @{ Cell = "B5"; Value = "Clicks (IOne)" }
)

$convEdits = @(
    @{ Cell = "B3"; Value = "Click Based Conversions (IOne)" },
    @{ Cell = "B4"; Value = "View Based Conversions (IOne)" }
)

# Sheets 1-4: CreativeDelivery_* (Impressions / Cost / Clicks rows)
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($edit in $deliveryEdits) {
        $ws.Range($edit.Cell).Value = $edit.Value
    }
}

# Sheets 5-8: CreativeConv_* (Click/View Based Conversions rows)
for ($i = 5; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($edit in $convEdits) {
        $ws.Range($edit.Cell).Value = $edit.Value
    }
}

# Restore the per-sheet selection state (cosmetic, matches the authored
# commit's saved view) and move the active tab from sheet 5 to sheet 8.
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Activate()
[void]$ws1.Range("C5").Select()

$ws2 = $wb.Worksheets.Item(2)
[void]$ws2.Activate()
[void]$ws2.Range("B3:B5").Select()

$ws3 = $wb.Worksheets.Item(3)
[void]$ws3.Activate()
[void]$ws3.Range("B3:B5").Select()

$ws4 = $wb.Worksheets.Item(4)
[void]$ws4.Activate()
[void]$ws4.Range("B3:B5").Select()

$ws5 = $wb.Worksheets.Item(5)
[void]$ws5.Activate()
[void]$ws5.Range("B3:B4").Select()

$ws6 = $wb.Worksheets.Item(6)
[void]$ws6.Activate()
[void]$ws6.Range("B3:B4").Select()

$ws7 = $wb.Worksheets.Item(7)
[void]$ws7.Activate()
[void]$ws7.Range("B3:B4").Select()

$ws8 = $wb.Worksheets.Item(8)
[void]$ws8.Activate()
[void]$ws8.Range("B3:B4").Select()
